# Auto-generated edit script: updates currentAveragePrice / profit columns
# (columns H-N) for specific leve rows across all 8 job sheets, per the
# scheduled price-refresh run. Values come from the upstream market-board
# snapshot; a few rows also gain/lose the optional M/N (profit) cell
# depending on whether NQ/HQ data is available for that refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 9
$ws.Range("H9").Value = 735.625
$ws.Range("I9").Value = 277.2
$ws.Range("J9").Value = 1499.6666
$ws.Range("K9").Value = 277.2
$ws.Range("L9").Value = 1499.6666
$ws.Range("M9").Value = -108.2
$ws.Range("N9").Value = -1837.6666

# Row 33
$ws.Range("H33").Value = 2644.5
$ws.Range("I33").Value = 2826
$ws.Range("J33").Value = 2100
$ws.Range("K33").Value = 2826
$ws.Range("L33").Value = 2100
$ws.Range("M33").Value = -2597
$ws.Range("N33").Value = -2558

# Row 88
$ws.Range("H88").Value = 17596802
$ws.Range("I88").Value = 47627336
$ws.Range("K88").Value = 47627336
$ws.Range("M88").Value = -47626930

# Row 91
$ws.Range("H91").Value = 17596802
$ws.Range("I91").Value = 47627336
$ws.Range("K91").Value = 47627336
$ws.Range("M91").Value = -47625932

# Row 106
$ws.Range("H106").Value = 1553.56
$ws.Range("I106").Value = 1098.5454
$ws.Range("J106").Value = 4890.3335
$ws.Range("K106").Value = 1098.5454
$ws.Range("L106").Value = 4890.3335
$ws.Range("M106").Value = -467.5454
$ws.Range("N106").Value = -6152.3335

# Row 141
$ws.Range("H141").Value = 1661.7878
$ws.Range("I141").Value = 1478.8064
$ws.Range("J141").Value = 4498
$ws.Range("K141").Value = 4436.4192
$ws.Range("L141").Value = 13494
$ws.Range("M141").Value = 743.5807999999997
$ws.Range("N141").Value = -23854


$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 2124164.5
$ws.Range("I32").Value = 2197380.5
$ws.Range("K32").Value = 2197380.5
$ws.Range("M32").Value = -2197093.5

# Row 61
$ws.Range("H61").Value = 6488.4185
$ws.Range("I61").Value = 3857.5557
$ws.Range("K61").Value = 3857.5557
$ws.Range("M61").Value = -3645.5557

# Row 63
$ws.Range("H63").Value = 6000.857
$ws.Range("I63").Value = 3401.2
$ws.Range("J63").Value = 12500
$ws.Range("K63").Value = 3401.2
$ws.Range("L63").Value = 12500
$ws.Range("M63").Value = -2715.2
$ws.Range("N63").Value = -13872

# Row 66
$ws.Range("H66").Value = 6000.857
$ws.Range("I66").Value = 3401.2
$ws.Range("J66").Value = 12500
$ws.Range("K66").Value = 17006
$ws.Range("L66").Value = 62500
$ws.Range("M66").Value = -13574
$ws.Range("N66").Value = -69364

# Row 74
$ws.Range("H74").Value = 27795.125
$ws.Range("I74").Value = 42180.24
$ws.Range("K74").Value = 42180.24
$ws.Range("M74").Value = -41306.24

# Row 77
$ws.Range("H77").Value = 27795.125
$ws.Range("I77").Value = 42180.24
$ws.Range("K77").Value = 210901.2
$ws.Range("M77").Value = -206533.2

# Row 88
$ws.Range("H88").Value = 6301.2
$ws.Range("I88").Value = 4753
$ws.Range("J88").Value = 7333.3335
$ws.Range("K88").Value = 4753
$ws.Range("L88").Value = 7333.3335
$ws.Range("M88").Value = -4347
$ws.Range("N88").Value = -8145.3335

# Row 91
$ws.Range("H91").Value = 6301.2
$ws.Range("I91").Value = 4753
$ws.Range("J91").Value = 7333.3335
$ws.Range("K91").Value = 4753
$ws.Range("L91").Value = 7333.3335
$ws.Range("M91").Value = -3349
$ws.Range("N91").Value = -10141.3335

# Row 122
$ws.Range("H122").Value = 3274.3
$ws.Range("I122").Value = 2128.3225
$ws.Range("K122").Value = 6384.967500000001
$ws.Range("M122").Value = -3934.967500000001

# Row 136
$ws.Range("H136").Value = 6488.4185
$ws.Range("I136").Value = 3857.5557
$ws.Range("K136").Value = 11572.6671
$ws.Range("M136").Value = -9022.667099999999


$ws = $wb.Worksheets.Item("BSM")

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 64
$ws.Range("H64").Value = 9260701
$ws.Range("I64").Value = 20834204
$ws.Range("J64").Value = 1898.4
$ws.Range("K64").Value = 20834204
$ws.Range("L64").Value = 1898.4
$ws.Range("M64").Value = -20833979
$ws.Range("N64").Value = -2348.4

# Row 67
$ws.Range("H67").Value = 9260701
$ws.Range("I67").Value = 20834204
$ws.Range("J67").Value = 1898.4
$ws.Range("K67").Value = 20834204
$ws.Range("L67").Value = 1898.4
$ws.Range("M67").Value = -20833424
$ws.Range("N67").Value = -3458.4

# Row 102
$ws.Range("H102").Value = 4749
$ws.Range("I102").Value = 4749
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4749
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1504
$ws.Range("N102").ClearContents()

# Row 105
$ws.Range("H105").Value = 3465.4644
$ws.Range("I105").Value = 1816.2727
$ws.Range("K105").Value = 1816.2727
$ws.Range("M105").Value = -69.27269999999999


$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 6754.932
$ws.Range("I31").Value = 3609.879
$ws.Range("K31").Value = 3609.879
$ws.Range("M31").Value = -3314.879

# Row 34
$ws.Range("H34").Value = 6754.932
$ws.Range("I34").Value = 3609.879
$ws.Range("K34").Value = 3609.879
$ws.Range("M34").Value = -3407.879

# Row 122
$ws.Range("H122").Value = 2410.2068
$ws.Range("I122").Value = 1507.95
$ws.Range("K122").Value = 4523.85
$ws.Range("M122").Value = -2073.85

# Row 132
$ws.Range("H132").Value = 4454.377
$ws.Range("I132").Value = 2627
$ws.Range("J132").Value = 7935.095
$ws.Range("K132").Value = 7881
$ws.Range("L132").Value = 23805.285
$ws.Range("M132").Value = -5351
$ws.Range("N132").Value = -28865.285


$ws = $wb.Worksheets.Item("CUL")

# Row 33
$ws.Range("H33").Value = 30303304
$ws.Range("I33").Value = 66666748
$ws.Range("J33").Value = 433.83334
$ws.Range("K33").Value = 400000488
$ws.Range("L33").Value = 2603.00004
$ws.Range("M33").Value = -400000205
$ws.Range("N33").Value = -3169.00004

# Row 113
$ws.Range("H113").Value = 3933.25
$ws.Range("J113").Value = 5199.8184
$ws.Range("L113").Value = 15599.4552
$ws.Range("N113").Value = -19939.4552

# Row 132
$ws.Range("H132").Value = 6560.303
$ws.Range("I132").Value = 3286.6667
$ws.Range("J132").Value = 9288.333000000001
$ws.Range("K132").Value = 29580.0003
$ws.Range("L132").Value = 83594.997
$ws.Range("M132").Value = -27050.0003
$ws.Range("N132").Value = -88654.997

# Row 140
$ws.Range("H140").Value = 107910.9
$ws.Range("I140").Value = 154947.39
$ws.Range("J140").Value = 5998.5
$ws.Range("K140").Value = 464842.17
$ws.Range("L140").Value = 17995.5
$ws.Range("M140").Value = -459662.17
$ws.Range("N140").Value = -28355.5


$ws = $wb.Worksheets.Item("GSM")

# Row 126
$ws.Range("H126").Value = 8012
$ws.Range("I126").Value = 8012
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 24036
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -21566
$ws.Range("N126").ClearContents()


$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Range("H46").Value = 1727636.8
$ws.Range("I46").Value = 4311661
$ws.Range("J46").Value = 4953.8335
$ws.Range("K46").Value = 4311661
$ws.Range("L46").Value = 4953.8335
$ws.Range("M46").Value = -4311473
$ws.Range("N46").Value = -5329.8335

# Row 132
$ws.Range("H132").Value = 7817615.5
$ws.Range("J132").Value = 8877.593000000001
$ws.Range("L132").Value = 26632.779
$ws.Range("N132").Value = -31692.779


$ws = $wb.Worksheets.Item("WVR")

# Row 132
$ws.Range("H132").Value = 15387650
$ws.Range("I132").Value = 18184696
$ws.Range("K132").Value = 54554088
$ws.Range("M132").Value = -54551558

